$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "Status code on 15-Oct-2020-12-31-06"
$ws.Range("F1").Value = "Content type on 15-Oct-2020-12-31-06"
$ws.Range("G1").Value = "Response body on 15-Oct-2020-12-31-06"
